$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("set check")
$ws.Range("A35").Value = "num"
